$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 5.147176742553711
$ws.Range("B1").Value = 5.374391078948975
$ws.Range("C1").Value = 5.378048896789551
$ws.Range("D1").Value = 9.033856391906738
$ws.Range("E1").Value = 7.45313835144043
